# Daily attendance processing - 2025-12-30 07:37:09
# Normalizes the "Recorded By" (column G) comma-separated recorder lists:
# the first recorder in the list is rotated to the end of the list for a
# known set of recorder-combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact before -> after mapping observed for the "Recorded By" column.
$map = @{
    "System, backup@backdoor.com, system" = "backup@backdoor.com, system, System";
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
